# The paragraph currently reads (across 4 runs):
#   "Versi" | "on" | " 2" | "."   (bookmark _GoBack sits between " 2" and ".")
# and must become (across 2 runs):
#   "Version" | " 1."              (bookmark _GoBack now sits after " 1.")

$d = $word.ActiveDocument

# 1) Merge the "Versi" + "on" runs into a single "Version" run.
$rngVersion = $d.Range(0, 7)
$rngVersion.Find.Execute("Version", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Version", 2)

# 2) Change the "2" in the " 2" run to "1." (keeps this edit inside that run,
#    so it doesn't cross the _GoBack bookmark and delete it).
$rngNumber = $d.Range(7, 9)
$rngNumber.Find.Execute("2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.", 2)

# 3) Remove the now-redundant trailing "." run (the period that used to
#    follow the bookmark) - again without crossing the bookmark itself.
$rngTail = $d.Range(10, 11)
$rngTail.Find.Execute(".", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)
